$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)
try { Write-Output ($lo.TableObject | Out-String) } catch { Write-Output "err1: $_" }
try { Get-Member -InputObject $lo.TableObject | Out-String | Write-Output } catch { Write-Output "err2: $_" }
